$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the reference value in row 3 (E3): RE529419 -> R523409
$ws.Range("E3").Value = "R523409"

# Update the confidence value in row 3 (I3)
$ws.Range("I3").Value = 89.906999999999996

# Recalculate so dependent formulas (I12 AVERAGEIF) pick up the new value
$excel.Calculate()

# Update the active selection to I2:I10 with active cell I2
$ws.Activate()
$ws.Range("I2:I10").Select()
